# Delete row 593 ("「無限の宇宙を旅した光」...") entirely.
# This shifts all subsequent rows up by one, matching the diff
# (old row 594 becomes new row 593, ... old row 654 becomes new row 653).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(593).Delete()
